$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 98, shifting the
# existing rows 98-177 down to 100-179 (weekly price update: a new
# "Cilantro" observation for 2023-02-24 is added while keeping all the
# previously recorded observations intact).
$ws.Rows("98:99").Insert()

# New row 98 - "Primera" quality observation for 2023-02-24
$ws.Cells.Item(98, 1).Value = 7
$ws.Cells.Item(98, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(98, 3).Value = "Ñuble"
$ws.Cells.Item(98, 4).Value = "2023-02-24"
$ws.Cells.Item(98, 5).Value = 16
$ws.Cells.Item(98, 6).Value = 100112040
$ws.Cells.Item(98, 7).Value = "Cilantro"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 500
$ws.Cells.Item(98, 11).Value = 700
$ws.Cells.Item(98, 12).Value = 700
$ws.Cells.Item(98, 13).Value = 700
$ws.Cells.Item(98, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(98, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(98, 16).Value = 700
$ws.Cells.Item(98, 17).Value = 1
$ws.Cells.Item(98, 18).Value = "Hortaliza"

# New row 99 - "Segunda" quality observation for 2023-02-24
$ws.Cells.Item(99, 1).Value = 7
$ws.Cells.Item(99, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(99, 3).Value = "Ñuble"
$ws.Cells.Item(99, 4).Value = "2023-02-24"
$ws.Cells.Item(99, 5).Value = 16
$ws.Cells.Item(99, 6).Value = 100112040
$ws.Cells.Item(99, 7).Value = "Cilantro"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Segunda"
$ws.Cells.Item(99, 10).Value = 300
$ws.Cells.Item(99, 11).Value = 600
$ws.Cells.Item(99, 12).Value = 600
$ws.Cells.Item(99, 13).Value = 600
$ws.Cells.Item(99, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(99, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(99, 16).Value = 600
$ws.Cells.Item(99, 17).Value = 1
$ws.Cells.Item(99, 18).Value = "Hortaliza"

# Make sure the date cells use the same date/time display format as the
# rest of the "Fecha" column.
$ws.Range("D98:D99").NumberFormat = "YYYY-MM-DD HH:MM:SS"
